$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.7310083333333334
$ws.Range("H2").Value = 2.193025
$ws.Range("I2").Value = 0.01673731480740535
$ws.Range("J2").Value = 0.01673731480740535
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3515896666666666
$ws.Range("N2").Value = 1.054769
$ws.Range("O2").Value = 0.03702657605490316
$ws.Range("P2").Value = 0.03702657605490316
$ws.Range("Q2").Value = 0.2570149762472222
$ws.Range("R2").Value = 2.313134786225
$ws.Range("S2").Value = 0.000619725459671251
$ws.Range("T2").Value = 0.000619725459671251

$ws.Range("G3").Value = 0.7310083333333334
$ws.Range("H3").Value = 2.193025
$ws.Range("I3").Value = 0.01673731480740535
$ws.Range("J3").Value = 0.01673731480740535
$ws.Range("O3").Value = 0.1974721703648871
$ws.Range("P3").Value = 0.1974721703648871
$ws.Range("Q3").Value = 1.370726396644444
$ws.Range("R3").Value = 12.3365375698
$ws.Range("S3").Value = 0.003305153881098697
$ws.Range("T3").Value = 0.003305153881098697

$ws.Range("G4").Value = 0.7310083333333334
$ws.Range("H4").Value = 2.193025
$ws.Range("I4").Value = 0.01673731480740535
$ws.Range("J4").Value = 0.01673731480740535
$ws.Range("M4").Value = 7.268896000000001
$ws.Range("N4").Value = 21.806688
$ws.Range("O4").Value = 0.7655012535802097
$ws.Range("P4").Value = 0.7655012535802097
$ws.Range("Q4").Value = 5.313623550133334
$ws.Range("R4").Value = 47.8226119512
$ws.Range("S4").Value = 0.0128124354666354
$ws.Range("T4").Value = 0.0128124354666354

$ws.Range("I5").Value = 0.8536212576586365
$ws.Range("J5").Value = 0.8536212576586365
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.3515896666666666
$ws.Range("N5").Value = 1.054769
$ws.Range("O5").Value = 0.03702657605490316
$ws.Range("P5").Value = 0.03702657605490316
$ws.Range("Q5").Value = 13.10804330239333
$ws.Range("R5").Value = 117.97238972154
$ws.Range("S5").Value = 0.03160667241877959
$ws.Range("T5").Value = 0.03160667241877959

$ws.Range("I6").Value = 0.8536212576586365
$ws.Range("J6").Value = 0.8536212576586365
$ws.Range("O6").Value = 0.1974721703648871
$ws.Range("P6").Value = 0.1974721703648871
$ws.Range("Q6").Value = 69.90853694714664
$ws.Range("R6").Value = 629.1768325243198
$ws.Range("S6").Value = 0.1685664424194555
$ws.Range("T6").Value = 0.1685664424194555

$ws.Range("I7").Value = 0.8536212576586365
$ws.Range("J7").Value = 0.8536212576586365
$ws.Range("M7").Value = 7.268896000000001
$ws.Range("N7").Value = 21.806688
$ws.Range("O7").Value = 0.7655012535802097
$ws.Range("P7").Value = 0.7655012535802097
$ws.Range("Q7").Value = 271.00057982912
$ws.Range("R7").Value = 2439.00521846208
$ws.Range("S7").Value = 0.6534481428204014
$ws.Range("T7").Value = 0.6534481428204014

$ws.Range("G8").Value = 5.662136666666666
$ws.Range("H8").Value = 16.98641
$ws.Range("I8").Value = 0.129641427533958
$ws.Range("J8").Value = 0.129641427533958
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.3515896666666666
$ws.Range("N8").Value = 1.054769
$ws.Range("O8").Value = 0.03702657605490316
$ws.Range("P8").Value = 0.03702657605490316
$ws.Range("Q8").Value = 1.990748743254444
$ws.Range("R8").Value = 17.91673868929
$ws.Range("S8").Value = 0.004800178176452313
$ws.Range("T8").Value = 0.004800178176452313

$ws.Range("G9").Value = 5.662136666666666
$ws.Range("H9").Value = 16.98641
$ws.Range("I9").Value = 0.129641427533958
$ws.Range("J9").Value = 0.129641427533958
$ws.Range("O9").Value = 0.1974721703648871
$ws.Range("P9").Value = 0.1974721703648871
$ws.Range("Q9").Value = 10.61717060736889
$ws.Range("R9").Value = 95.55453546631999
$ws.Range("S9").Value = 0.02560057406433292
$ws.Range("T9").Value = 0.02560057406433292

$ws.Range("G10").Value = 5.662136666666666
$ws.Range("H10").Value = 16.98641
$ws.Range("I10").Value = 0.129641427533958
$ws.Range("J10").Value = 0.129641427533958
$ws.Range("M10").Value = 7.268896000000001
$ws.Range("N10").Value = 21.806688
$ws.Range("O10").Value = 0.7655012535802097
$ws.Range("P10").Value = 0.7655012535802097
$ws.Range("Q10").Value = 41.15748256778667
$ws.Range("R10").Value = 370.41734311008
$ws.Range("S10").Value = 0.09924067529317278
$ws.Range("T10").Value = 0.09924067529317278

